$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.907249101897178
$ws.Range("BO1").Value = 0.70783710240559872
$ws.Range("BP1").Value = 0.81700234240870917
$ws.Range("D2").Value = 0.81374293228114425
$ws.Range("BP2").Value = 0.72348244641236459
$ws.Range("A3").Value = 0.99592686518003193
$ws.Range("B3").Value = 0.9433722397428933
$ws.Range("E3").Value = 0.87300552271672371
$ws.Range("AL3").Value = 0.91816350672036129
$ws.Range("E4").Value = 0.70541149368193823
$ws.Range("F5").Value = 0.97773723464188156
$ws.Range("E7").Value = 0.93496117610854501
$ws.Range("F7").Value = 0.93500546928525408
$ws.Range("I7").Value = 0.97032776308452706
$ws.Range("F8").Value = 0.81062779818158481
$ws.Range("I8").Value = 0.91627418835532626
$ws.Range("M8").Value = 0.96633008170315904
$ws.Range("H10").Value = 0.97560291787987397
$ws.Range("AL10").Value = 0.949263574308157
$ws.Range("AV10").Value = 0.894342148930934
$ws.Range("I11").Value = 0.91579753249569706
$ws.Range("J12").Value = 0.56660046914396145
$ws.Range("K12").Value = 0.8163238678166661
$ws.Range("N12").Value = 0.63103805208536234
$ws.Range("L13").Value = 0.78509489865679949
$ws.Range("S13").Value = 0.57297479897250381
$ws.Range("P14").Value = 0.89342612919605213
$ws.Range("AT14").Value = 0.82390487419862013
$ws.Range("M15").Value = 0.72754128862493928
$ws.Range("N15").Value = 0.91492819619526411
$ws.Range("Q16").Value = 0.85884516057155735
$ws.Range("R16").Value = 0.97009114076574465
$ws.Range("O17").Value = 0.67115520240967608
$ws.Range("R17").Value = 0.59182931748200218
$ws.Range("S18").Value = 0.58569493938210615
$ws.Range("T18").Value = 0.7571582755156947
$ws.Range("Q19").Value = 0.97018862636767533
$ws.Range("T19").Value = 0.97109998025900102
$ws.Range("V20").Value = 0.63584935832887302
$ws.Range("BE20").Value = 0.6152679011839981
$ws.Range("BN21").Value = 0.67698532858814664
$ws.Range("W22").Value = 0.83219544253872879
$ws.Range("X22").Value = 0.96485132648864402
$ws.Range("X23").Value = 0.90924598419404568
$ws.Range("Y23").Value = 0.84783032489727739
$ws.Range("AG24").Value = 0.59282817220643103
$ws.Range("X25").Value = 0.98593739856807505
$ws.Range("X26").Value = 0.96871769126933804
$ws.Range("K27").Value = 0.71899100302122188
$ws.Range("Z27").Value = 0.79530812216265034
$ws.Range("AA28").Value = 0.59449466003394913
$ws.Range("AD28").Value = 0.7034320434522896
$ws.Range("I29").Value = 0.95511805578844244
$ws.Range("AB29").Value = 0.95281693198066231
$ws.Range("AE29").Value = 0.69215615992444346
$ws.Range("AC30").Value = 0.99874622025931592
$ws.Range("AR30").Value = 0.56185029632050232
$ws.Range("AD31").Value = 0.8997382865265291
$ws.Range("AG31").Value = 0.57165702905517102
$ws.Range("AG32").Value = 0.90299667573422215
$ws.Range("AS32").Value = 0.85013735780669575
$ws.Range("G33").Value = 0.93095935363186155
$ws.Range("Y34").Value = 0.66987650024962453
$ws.Range("AF34").Value = 0.83555334788572655
$ws.Range("AG34").Value = 0.99756527483231983
$ws.Range("AH35").Value = 0.64341143155171876
$ws.Range("AJ35").Value = 0.80056474361439123
$ws.Range("AK35").Value = 0.89364355834418707
$ws.Range("BB35").Value = 0.87772536989475725
$ws.Range("U36").Value = 0.77159149606233202
$ws.Range("AL36").Value = 0.85462026933984592
$ws.Range("AY36").Value = 0.78540168418077183
$ws.Range("BA36").Value = 0.77686893565481663
$ws.Range("AL37").Value = 0.73206463625203755
$ws.Range("AA38").Value = 0.72774679758229754
$ws.Range("AM38").Value = 0.86622597584481509
$ws.Range("AK39").Value = 0.93429251579828398
$ws.Range("AO39").Value = 0.83549460857778901
$ws.Range("AU39").Value = 0.95848995548874361
$ws.Range("AL40").Value = 0.55134606420257626
$ws.Range("AZ40").Value = 0.841051956583586
$ws.Range("B41").Value = 0.87448648584875244
$ws.Range("AN42").Value = 0.97106467835067733
$ws.Range("AR42").Value = 0.86072765404364016
$ws.Range("AZ42").Value = 0.86251278979999113
$ws.Range("AO43").Value = 0.87241117003610169
$ws.Range("BE43").Value = 0.89011414929267518
$ws.Range("D44").Value = 0.90651763055132939
$ws.Range("AS44").Value = 0.92672905157828478
$ws.Range("O45").Value = 0.77185302700497105
$ws.Range("AT45").Value = 0.86312919362545681
$ws.Range("AE46").Value = 0.86728888148813699
$ws.Range("AR46").Value = 0.7221068336982015
$ws.Range("AT47").Value = 0.94732053775688196
$ws.Range("AW47").Value = 0.92356002250940228
$ws.Range("AK48").Value = 0.95727605584116437
$ws.Range("AU48").Value = 0.8618209408473152
$ws.Range("Y49").Value = 0.83055552273050526
$ws.Range("AX49").Value = 0.81073373903005919
$ws.Range("AZ50").Value = 0.77485620911080244
$ws.Range("AW51").Value = 0.72992205753812245
$ws.Range("AX51").Value = 0.91174014302865802
$ws.Range("AY53").Value = 0.87674552354277258
$ws.Range("AZ53").Value = 0.61161778462780703
$ws.Range("BB53").Value = 0.90493438840691132
$ws.Range("AZ54").Value = 0.85263287532061027
$ws.Range("BC54").Value = 0.81392109310040939
$ws.Range("Z55").Value = 0.87293286310248908
$ws.Range("BA55").Value = 0.98077857393607182
$ws.Range("BD55").Value = 0.85881267739249156
$ws.Range("BE55").Value = 0.93142683400912807
$ws.Range("AH56").Value = 0.79569616829044865
$ws.Range("BF56").Value = 0.98319150509419984
$ws.Range("Q57").Value = 0.69632467410636179
$ws.Range("BD57").Value = 0.72766357178631713
$ws.Range("Q58").Value = 0.84592993098418434
$ws.Range("BE58").Value = 0.97702100964061356
$ws.Range("BF59").Value = 0.93151777975834293
$ws.Range("BH59").Value = 0.93967656159412383
$ws.Range("AV60").Value = 0.65130002130584763
$ws.Range("AQ61").Value = 0.85859312558058809
$ws.Range("BG61").Value = 0.94033512597897784
$ws.Range("BK61").Value = 0.79180447796382447
$ws.Range("U62").Value = 0.96151115291120093
$ws.Range("BH62").Value = 0.78563152829766969
$ws.Range("BI62").Value = 0.96057032026408118
$ws.Range("BL62").Value = 0.58835277590924351
$ws.Range("BM63").Value = 0.86929706703337062
$ws.Range("BK64").Value = 0.90142056323801456
$ws.Range("BM64").Value = 0.56267829348144138
$ws.Range("D65").Value = 0.66731901521137327
$ws.Range("BL66").Value = 0.95460872462848112
$ws.Range("BM66").Value = 0.88738052193590278
$ws.Range("BN67").Value = 0.63585638059374716
$ws.Range("BP67").Value = 0.93369307475193741
$ws.Range("BN68").Value = 0.63387976368469201
